$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New timeline entries (Chapter 6 content)
$newRows = @(
    @{ Row = 42; Date = -50000; Text = "Austric people come to India." },
    @{ Row = 43; Date = -30000; Text = "Dravidians come to India." },
    @{ Row = 44; Date = -4000;  Text = "Brahui (Dravidian) speakers come to India." },
    @{ Row = 45; Date = -500;   Text = "Start of the Middle Indo-Aryan, covering Prakrit, Pali, and Apabhramsha languages." },
    @{ Row = 46; Date = 1800;   Text = "Oral legends and traditions of Munda and Tibeto-Burman languages recorded by Christian missionaries." }
)

foreach ($entry in $newRows) {
    $ws.Cells.Item($entry.Row, 1).Value = $entry.Date
    $ws.Cells.Item($entry.Row, 3).Value = $entry.Text
}

# Update selection / active cell to mirror the new end of data
$ws.Application.Goto($ws.Range("C47"))
$ws.Range("C47").Select()
